# Update "想去人数" (F column) figures across sheets to match the
# latest generated output (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Cells.Item(2, 6).Value = 7658
$wsExpo.Cells.Item(6, 6).Value = 4327
$wsExpo.Cells.Item(7, 6).Value = 329
$wsExpo.Cells.Item(8, 6).Value = 605
$wsExpo.Cells.Item(10, 6).Value = 682

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(2, 6).Value = 51

# Sheet "全部类型" (all types, combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(2, 6).Value = 7658
$wsAll.Cells.Item(3, 6).Value = 51
$wsAll.Cells.Item(7, 6).Value = 4327
$wsAll.Cells.Item(8, 6).Value = 329
$wsAll.Cells.Item(9, 6).Value = 605
$wsAll.Cells.Item(11, 6).Value = 682
